# Generate Report for Handback
# Applies the "handback" report-generation edit to the localization-status workbook:
#  - updates the status text from "In Translation" to "Handed back: in sync with en-US"
#    everywhere it appears (Overview + per-locale sheets),
#  - fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#    columns on the zh-cn and de-de sheets for both rows, including hyperlinks on the
#    newly-populated "Latest Target File" cells,
#  - widens the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusOld = "In Translation"
$statusNew = "Handed back: in sync with en-US"

$mdFile1       = "1c7f37c1-7045-4405-8dd2-26a24c42bf57.md"
$mdFile2       = "9f3f175d-85e7-4c45-a47f-9db670de3280.md"
$url1          = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a8dd067c3a93c45d7b7ada5cff8ec6d693ee58a/e2e/1c7f37c1-7045-4405-8dd2-26a24c42bf57.md"
$url2          = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a8dd067c3a93c45d7b7ada5cff8ec6d693ee58a/e2e/9f3f175d-85e7-4c45-a47f-9db670de3280.md"

$zhHandback1   = "1c7f37c1-7045-4405-8dd2-26a24c42bf57.c6491033f6bf62556dc5f1641f2a2c79295dbbd6.zh-cn.xlf"
$zhHandback2   = "9f3f175d-85e7-4c45-a47f-9db670de3280.39b62c6f4dc079441bc25814fe77854fad88edd8.zh-cn.xlf"
$deHandback1   = "1c7f37c1-7045-4405-8dd2-26a24c42bf57.c6491033f6bf62556dc5f1641f2a2c79295dbbd6.de-de.xlf"
$deHandback2   = "9f3f175d-85e7-4c45-a47f-9db670de3280.39b62c6f4dc079441bc25814fe77854fad88edd8.de-de.xlf"

$zhHandbackDateTime = "2016-09-04 08:27:17"
$deHandbackDateTime = "2016-09-04 08:27:25"

# ---------------------------------------------------------------------------
# Overview sheet: refresh the status text shown for both rows
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $statusNew
$wsOverview.Range("F2").Value2 = $statusNew
$wsOverview.Range("E3").Value2 = $statusNew
$wsOverview.Range("F3").Value2 = $statusNew

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# Helper locale sheets: zh-cn and de-de each hold the same two rows of data,
# only the handback file names / datetime differ.
# ---------------------------------------------------------------------------
function Update-LocaleSheet($ws, $handback1, $handback2, $handbackDateTime) {

    # Status column (C)
    $ws.Range("C2").Value2 = $statusNew
    $ws.Range("C3").Value2 = $statusNew

    # Latest Handback File (J) and Latest Handback DateTime (K)
    $ws.Range("J2").Value2 = $handback1
    $ws.Range("K2").Value2 = $handbackDateTime
    $ws.Range("J3").Value2 = $handback2
    $ws.Range("K3").Value2 = $handbackDateTime

    # Latest Target File (I) values
    $ws.Range("I2").Value2 = $mdFile1
    $ws.Range("I3").Value2 = $mdFile2

    # Rebuild hyperlinks in row order (A2, I2, A3, I3) so the relationship ids
    # line up the same way Excel produces them when regenerating the report.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $url1, "", "", $mdFile1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $url1, "", "", $mdFile1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $url2, "", "", $mdFile2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $url2, "", "", $mdFile2)

    # Column widths: C, I and J now hold long text like the other file-name columns
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LocaleSheet $wsZhCn $zhHandback1 $zhHandback2 $zhHandbackDateTime

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LocaleSheet $wsDeDe $deHandback1 $deHandback2 $deHandbackDateTime
